$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-28 Friday" "2025-11-29 Saturday"

Replace-Text "805×8=6440" "589×2=1178"
Replace-Text "389×7=2723" "615×4=2460"
Replace-Text "197×8=1576" "460×9=4140"
Replace-Text "909×6=5454" "770×5=3850"
Replace-Text "163×8=1304" "213×6=1278"

Replace-Text "855×3=2565" "668×7=4676"
Replace-Text "278×4=1112" "629×6=3774"
Replace-Text "388×3=1164" "224×6=1344"
Replace-Text "879×7=6153" "811×3=2433"
Replace-Text "445×3=1335" "160×4=640"

Replace-Text "884×6=5304" "378×8=3024"
Replace-Text "573×8=4584" "185×4=740"
Replace-Text "366×3=1098" "936×2=1872"
Replace-Text "468×8=3744" "665×7=4655"
Replace-Text "897×2=1794" "838×2=1676"

Replace-Text "876×6=5256" "635×9=5715"
Replace-Text "219×5=1095" "720×9=6480"
Replace-Text "899×9=8091" "751×5=3755"
Replace-Text "266×5=1330" "240×8=1920"
Replace-Text "734×9=6606" "248×8=1984"

Replace-Text "881×8=7048" "334×2=668"
Replace-Text "361×9=3249" "970×4=3880"
Replace-Text "783×8=6264" "911×2=1822"
Replace-Text "363×8=2904" "813×9=7317"
Replace-Text "767×4=3068" "148×2=296"
